# Weekly update: a new reporting week (2021-11-09, serial 44509) is added
# for "Femacal de La Calera - Cebolla", inserted as 4 new rows at the top
# of this market's block (rows 760-763), pushing the existing history
# down by 4 rows (old row 760 -> 764, ..., old row 794 -> 798).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above row 760 (shifts 760:794 down to 764:798)
$ws.Range("A760:A763").EntireRow.Insert()

$rows = @(760, 761, 762, 763)

$A = @(3, 3, 3, 3)
$B = @("Femacal de La Calera", "Femacal de La Calera", "Femacal de La Calera", "Femacal de La Calera")
$C = @("Coquimbo", "Coquimbo", "Coquimbo", "Coquimbo")
$D = @(44509, 44509, 44509, 44509)
$E = @(5, 5, 5, 5)
$F = @(100112004, 100112004, 100112004, 100112004)
$G = @("Cebolla", "Cebolla", "Cebolla", "Cebolla")
$H = @("Morada(o)", "Sin especificar", "Sin especificar", "Sin especificar")
$I = @("1a nueva(o)", "1a (guarda)", "1a nueva(o)", "2a nueva(o)")
$J = @(80, 130, 14700, 7500)
$K = @(6000, 4500, 2300, 1500)
$L = @(6000, 4700, 2500, 1500)
$M = @(6000, 4592, 2398, 1500)
$N = @("`$/malla 18 kilos", "`$/malla 18 kilos", "`$/paquete 20 unidades (volumen en unidades)", "`$/paquete 20 unidades (volumen en unidades)")
$O = @("Región de Arica y Parinacota", "Provincia de Quillota", "Provincia de Quillota", "Provincia de Quillota")
$P = @(333, 255, 120, 75)
$Q = @(18, 18, 20, 20)
$R = @("Hortaliza", "Hortaliza", "Hortaliza", "Hortaliza")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $A[$i]
    $ws.Cells.Item($r, 2).Value = $B[$i]
    $ws.Cells.Item($r, 3).Value = $C[$i]
    $ws.Cells.Item($r, 4).Value = $D[$i]
    $ws.Cells.Item($r, 5).Value = $E[$i]
    $ws.Cells.Item($r, 6).Value = $F[$i]
    $ws.Cells.Item($r, 7).Value = $G[$i]
    $ws.Cells.Item($r, 8).Value = $H[$i]
    $ws.Cells.Item($r, 9).Value = $I[$i]
    $ws.Cells.Item($r, 10).Value = $J[$i]
    $ws.Cells.Item($r, 11).Value = $K[$i]
    $ws.Cells.Item($r, 12).Value = $L[$i]
    $ws.Cells.Item($r, 13).Value = $M[$i]
    $ws.Cells.Item($r, 14).Value = $N[$i]
    $ws.Cells.Item($r, 15).Value = $O[$i]
    $ws.Cells.Item($r, 16).Value = $P[$i]
    $ws.Cells.Item($r, 17).Value = $Q[$i]
    $ws.Cells.Item($r, 18).Value = $R[$i]
}
